# Add "Ant Colony Optimization" (ACO) into the Swarm-based group, ahead of
# "Artificial Bee Colony" (ABC), re-dating ACO to 2006, and bump the STT
# (sequence number) of every row that follows by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: becomes the "Ant Colony Optimization" entry (was ABC) -------
$ws.Range("C13").Value = "Ant Colony Optimization"
$ws.Range("D13").Value = "ACO"
$ws.Range("E13").Value = 2006
$ws.Range("F13").Value = "original"
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = "medium"

# --- Row 14: becomes the "Artificial Bee Colony" entry (was ACO) ---------
$ws.Range("C14").Value = "Artificial Bee Colony"
$ws.Range("D14").Value = "ABC"
$ws.Range("E14").Value = 2007
$ws.Range("F14").Value = "changed"
$ws.Range("K14").Value = 8
$ws.Range("L14").Value = "easy"

# --- Rows 15-37: the STT (B column) numbering shifts down by one ---------
for ($r = 15; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $cell.Value + 1
}

# --- View state: scroll so row 5 is at the top, select C14 ---------------
$ws.Activate()
$excel.Goto($ws.Range("A5"), $true)
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
